$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-25: column C originally held "Karate-Coed" / "Karate-Boys" / "Karate-Girls".
# New behaviour: column B encodes the coed/boys/girls variant, column C is just "Karate".
$boysRows = @(6, 17, 23)
$girlsRows = @(7, 18, 24)

for ($r = 2; $r -le 25; $r++) {
    if ($boysRows -contains $r) {
        $ws.Cells.Item($r, 2).Value = "sports_club_boys"
    }
    elseif ($girlsRows -contains $r) {
        $ws.Cells.Item($r, 2).Value = "sports_club_girls"
    }
    else {
        $ws.Cells.Item($r, 2).Value = "sports_club_coed"
    }

    $ws.Cells.Item($r, 3).Value = "Karate"
}
